# applied date normalization func
#
# The extraction pipeline's date-normalization pass turns bare
# "study start/end year" values into fully-qualified calendar dates:
#   - end-of-timeline years  -> "December 31, <year>"
#   - start-of-timeline years -> "January 1, <year>"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 (Lebanon cohort): study end year 2011 -> Dec 31, 2011
$ws.Range("U3").Value = "December 31, 2011"
# Row 4 (neuromuscular referral cohort): study end year 2019 -> Dec 31, 2019
$ws.Range("U4").Value = "December 31, 2019"
# Row 3: study start year 2010 -> Jan 1, 2010
$ws.Range("W3").Value = "January 1, 2010"
# Row 4: study start year 1999 -> Jan 1, 1999
$ws.Range("W4").Value = "January 1, 1999"

# The longer date strings no longer fit the old best-fit column widths;
# widen columns U ("study end year") and W ("study start year") to fit.
$ws.Columns.Item(21).ColumnWidth = 15.1665
$ws.Columns.Item(23).ColumnWidth = 12.333

# Reflect the reviewer's final on-screen state: scrolled right so the
# newly-edited study-year columns are visible, with W14 selected.
$ws.Range("W14").Select()
